# Fruta / hortaliza, semanal
# A new week of data (2021-12-28) is inserted at the top of the data block
# (right after the header-adjacent existing rows), pushing every existing
# data row down by two. The two new rows re-use the "Primera"/"Segunda"
# Peru price pair that was previously recorded for 2021-12-23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 20; everything from the old
# row 20 onward shifts down by two rows (old row 20 -> new row 22, etc.),
# and the sheet's used range grows from A1:T137 to A1:T139.
$ws.Range("A20:A21").EntireRow.Insert()

# New row 20: "Primera" quality entry for the week of 2021-12-28.
$ws.Cells.Item(20, 1).Value2 = 4
$ws.Cells.Item(20, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(20, 3).Value2 = "Los Lagos"
$ws.Cells.Item(20, 4).Value2 = 44558
$ws.Cells.Item(20, 5).Value2 = 10
$ws.Cells.Item(20, 6).Value2 = "Fruta"
$ws.Cells.Item(20, 7).Value2 = 100108
$ws.Cells.Item(20, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(20, 9).Value2 = 100108002
$ws.Cells.Item(20, 10).Value2 = "Mango"
$ws.Cells.Item(20, 11).Value2 = "Sin especificar"
$ws.Cells.Item(20, 12).Value2 = "Primera"
$ws.Cells.Item(20, 13).Value2 = 200
$ws.Cells.Item(20, 14).Value2 = 8500
$ws.Cells.Item(20, 15).Value2 = 9000
$ws.Cells.Item(20, 16).Value2 = 8750
$ws.Cells.Item(20, 17).Value2 = "$/bandeja 4 kilos"
$ws.Cells.Item(20, 18).Value2 = "Perú"
$ws.Cells.Item(20, 19).Value2 = 2188
$ws.Cells.Item(20, 20).Value2 = 4

# New row 21: "Segunda" quality entry for the same week.
$ws.Cells.Item(21, 1).Value2 = 4
$ws.Cells.Item(21, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(21, 3).Value2 = "Los Lagos"
$ws.Cells.Item(21, 4).Value2 = 44558
$ws.Cells.Item(21, 5).Value2 = 10
$ws.Cells.Item(21, 6).Value2 = "Fruta"
$ws.Cells.Item(21, 7).Value2 = 100108
$ws.Cells.Item(21, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(21, 9).Value2 = 100108002
$ws.Cells.Item(21, 10).Value2 = "Mango"
$ws.Cells.Item(21, 11).Value2 = "Sin especificar"
$ws.Cells.Item(21, 12).Value2 = "Segunda"
$ws.Cells.Item(21, 13).Value2 = 60
$ws.Cells.Item(21, 14).Value2 = 6000
$ws.Cells.Item(21, 15).Value2 = 6000
$ws.Cells.Item(21, 16).Value2 = 6000
$ws.Cells.Item(21, 17).Value2 = "$/bandeja 4 kilos"
$ws.Cells.Item(21, 18).Value2 = "Perú"
$ws.Cells.Item(21, 19).Value2 = 1500
$ws.Cells.Item(21, 20).Value2 = 4
